$wb = $excel.ActiveWorkbook

# Add the new "UpdateActivity" worksheet after the last existing sheet (MoreAttendees).
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$newSheet.Name = "UpdateActivity"

# Header row (bold, centered) - mirrors the layout used on the other data sheets.
$newSheet.Range("A1").Value = "Subject"
$newSheet.Range("B1").Value = "IndustryGroup"
$newSheet.Range("C1").Value = "ProductType"
$newSheet.Range("D1").Value = "Description"
$newSheet.Range("E1").Value = "MeetingNotes"
$newSheet.Range("F1").Value = "ExtAttendee"
$newSheet.Range("G1").Value = "HLAttendee"
$newSheet.Range("A1:G1").Font.Bold = $true
$newSheet.Range("A1:G1").HorizontalAlignment = -4108

# Row 2 - primary attendee test data.
$newSheet.Range("A2").Value = "Updated Primary Attendee"
$newSheet.Range("B2").Value = "FIG - Financial Institutions"
$newSheet.Range("C2").Value = "Advisory"
$newSheet.Range("D2").Value = "Updated Test Description Primary Attendee"
$newSheet.Range("E2").Value = "Updated Notes"
$newSheet.Range("F2").Value = "Test James"
$newSheet.Range("G2").Value = "Amanda Donovan"

# Row 3 - non-primary attendee test data.
$newSheet.Range("A3").Value = "Updated Non Primary Attendee"
$newSheet.Range("B3").Value = "FIG - Financial Institutions"
$newSheet.Range("C3").Value = "Advisory"
$newSheet.Range("D3").Value = "Updated Test Description Non Primary Attendee"
$newSheet.Range("E3").Value = "Updated Notes"
$newSheet.Range("F3").Value = "Test James"
$newSheet.Range("G3").Value = "Amanda Donovan"

# Column widths (matching the widths Excel computed for this data when it was authored).
$newSheet.Columns.Item(1).ColumnWidth = 27.0533854166667
$newSheet.Columns.Item(2).ColumnWidth = 24.3854166666667
$newSheet.Columns.Item(3).ColumnWidth = 13.4986979166667
$newSheet.Columns.Item(4).ColumnWidth = 42.1666666666667
$newSheet.Columns.Item(5).ColumnWidth = 13.7213541666667
$newSheet.Columns.Item(6).ColumnWidth = 10.7213541666667
$newSheet.Columns.Item(7).ColumnWidth = 16.8307291666667

# The new sheet becomes the active / selected tab, with C19 selected (leftover from manual entry).
$newSheet.Activate()
$newSheet.Range("C19").Select() | Out-Null
